# Applies weekly price/fruit-vegetable data corrections to rows 2-26
# (Fruta / hortaliza, semanal) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 45044
$ws.Range("L2").Value2 = "Especial"
$ws.Range("M2").Value2 = 40
$ws.Range("N2").Value2 = 13000
$ws.Range("O2").Value2 = 13000
$ws.Range("P2").Value2 = 13000
$ws.Range("S2").Value2 = 722

# Row 3
$ws.Range("D3").Value2 = 45044
$ws.Range("L3").Value2 = "Primera"
$ws.Range("N3").Value2 = 12000
$ws.Range("O3").Value2 = 12000
$ws.Range("P3").Value2 = 12000
$ws.Range("S3").Value2 = 667

# Row 4
$ws.Range("D4").Value2 = 45043
$ws.Range("L4").Value2 = "Especial"
$ws.Range("N4").Value2 = 13000
$ws.Range("O4").Value2 = 13000
$ws.Range("P4").Value2 = 13000
$ws.Range("S4").Value2 = 722

# Row 5
$ws.Range("D5").Value2 = 45043
$ws.Range("L5").Value2 = "Primera"
$ws.Range("M5").Value2 = 50
$ws.Range("N5").Value2 = 12000
$ws.Range("O5").Value2 = 12000
$ws.Range("P5").Value2 = 12000
$ws.Range("S5").Value2 = 667

# Row 6
$ws.Range("D6").Value2 = 44699
$ws.Range("L6").Value2 = "Especial"
$ws.Range("M6").Value2 = 60
$ws.Range("N6").Value2 = 13000
$ws.Range("O6").Value2 = 13000
$ws.Range("P6").Value2 = 13000
$ws.Range("R6").Value2 = "Provincia de Curicó"
$ws.Range("S6").Value2 = 867

# Row 7
$ws.Range("D7").Value2 = 44699
$ws.Range("L7").Value2 = "Primera"
$ws.Range("M7").Value2 = 120
$ws.Range("N7").Value2 = 11000
$ws.Range("O7").Value2 = 12000
$ws.Range("P7").Value2 = 11500
$ws.Range("R7").Value2 = "Provincia de Curicó"
$ws.Range("S7").Value2 = 767

# Row 8
$ws.Range("D8").Value2 = 45050
$ws.Range("M8").Value2 = 50

# Row 9
$ws.Range("D9").Value2 = 45050

# Row 10
$ws.Range("D10").Value2 = 45071
$ws.Range("M10").Value2 = 40
$ws.Range("N10").Value2 = 12000
$ws.Range("O10").Value2 = 12000
$ws.Range("P10").Value2 = 12000
$ws.Range("S10").Value2 = 667

# Row 11
$ws.Range("D11").Value2 = 45071
$ws.Range("L11").Value2 = "Segunda"
$ws.Range("M11").Value2 = 40
$ws.Range("N11").Value2 = 10000
$ws.Range("O11").Value2 = 10000
$ws.Range("P11").Value2 = 10000
$ws.Range("S11").Value2 = 556

# Row 12
$ws.Range("L12").Value2 = "Especial"
$ws.Range("M12").Value2 = 50
$ws.Range("N12").Value2 = 13000
$ws.Range("O12").Value2 = 13000
$ws.Range("P12").Value2 = 13000
$ws.Range("S12").Value2 = 722

# Row 13
$ws.Range("D13").Value2 = 45049
$ws.Range("L13").Value2 = "Primera"
$ws.Range("N13").Value2 = 12000
$ws.Range("O13").Value2 = 12000
$ws.Range("P13").Value2 = 12000
$ws.Range("Q13").Value2 = "`$/caja 18 kilos empedrada"
$ws.Range("R13").Value2 = "Región de O'Higgins"
$ws.Range("S13").Value2 = 667
$ws.Range("T13").Value2 = 18

# Row 14
$ws.Range("D14").Value2 = 45069
$ws.Range("M14").Value2 = 60
$ws.Range("N14").Value2 = 12000
$ws.Range("P14").Value2 = 12000
$ws.Range("Q14").Value2 = "`$/caja 18 kilos empedrada"
$ws.Range("R14").Value2 = "Región de O'Higgins"
$ws.Range("S14").Value2 = 667
$ws.Range("T14").Value2 = 18

# Row 15
$ws.Range("D15").Value2 = 45069
$ws.Range("L15").Value2 = "Segunda"
$ws.Range("M15").Value2 = 40
$ws.Range("N15").Value2 = 10000
$ws.Range("O15").Value2 = 10000
$ws.Range("P15").Value2 = 10000
$ws.Range("S15").Value2 = 556

# Row 16
$ws.Range("D16").Value2 = 45020
$ws.Range("M16").Value2 = 60
$ws.Range("Q16").Value2 = "`$/caja 18 kilos granel"

# Row 17
$ws.Range("D17").Value2 = 45040

# Row 18
$ws.Range("D18").Value2 = 45040
$ws.Range("M18").Value2 = 40

# Row 19
$ws.Range("D19").Value2 = 45070
$ws.Range("M19").Value2 = 60
$ws.Range("N19").Value2 = 10000
$ws.Range("O19").Value2 = 10000
$ws.Range("P19").Value2 = 10000
$ws.Range("Q19").Value2 = "`$/caja 18 kilos empedrada"
$ws.Range("S19").Value2 = 556

# Row 20
$ws.Range("D20").Value2 = 45062

# Row 21
$ws.Range("D21").Value2 = 45062
$ws.Range("M21").Value2 = 50

# Row 22
$ws.Range("D22").Value2 = 45033
$ws.Range("M22").Value2 = 60

# Row 23
$ws.Range("D23").Value2 = 45033
$ws.Range("M23").Value2 = 80

# Row 24
$ws.Range("D24").Value2 = 45076
$ws.Range("M24").Value2 = 30
$ws.Range("Q24").Value2 = "`$/caja 15 kilos granel"
$ws.Range("S24").Value2 = 800
$ws.Range("T24").Value2 = 15

# Row 25
$ws.Range("D25").Value2 = 45076
$ws.Range("L25").Value2 = "Segunda"
$ws.Range("M25").Value2 = 30
$ws.Range("N25").Value2 = 10000
$ws.Range("O25").Value2 = 10000
$ws.Range("P25").Value2 = 10000
$ws.Range("Q25").Value2 = "`$/caja 15 kilos granel"
$ws.Range("S25").Value2 = 667
$ws.Range("T25").Value2 = 15

# Row 26
$ws.Range("D26").Value2 = 45021
$ws.Range("M26").Value2 = 50
$ws.Range("Q26").Value2 = "`$/caja 18 kilos granel"
